# Re-generate the "Общее время" (total time) column so minutes and
# seconds are always zero-padded to two digits (hours stay unpadded),
# e.g. "11 ч. 3 мин. 26 сек." -> "11 ч. 03 мин. 26 сек."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column I holds "Общее время" (total time), data starts on row 2.
$timeCol = 9

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $timeCol)
    $v = $cell.Value2

    if ($v -ne $null -and $v -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours = $matches[1]
        $minutes = $matches[2].PadLeft(2, '0')
        $seconds = $matches[3].PadLeft(2, '0')
        $newValue = "$hours ч. $minutes мин. $seconds сек."

        if ($newValue -ne $v) {
            $cell.Value2 = $newValue
        }
    }
}
